# B1--and-B2-PowerPoint.pptx edit
#
# 1) The single table on the deck (slide 5) switches from the custom
#    "Table_0" style defined in tableStyles.xml to the built-in table
#    style {6ED3D6C1-355D-42D4-8B84-029545941B43}.
#
# 2) The theme applied to the slide master/presentation swaps from the
#    "Integral" / "Red Violet" palette to the stock "Office Theme" /
#    "Office" palette - i.e. every theme color swatch is recolored to
#    the default Office values.

$p = $ppt.ActivePresentation

# --- 1) Retarget the table's style -----------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{6ED3D6C1-355D-42D4-8B84-029545941B43}")
        }
    }
}

# --- 2) Recolor the theme from "Integral" (Red Violet) to "Office Theme" (Office) ---
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1      -> 000000
$cs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477   # folHlink -> 954F72
